$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.378.72"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.935.64"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7456"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "245.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3171"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06982"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07996"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "1.937.61"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.352"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.28%  "
$ws.Range("D17").Value = "30.390.58"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "252.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007925"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.727"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "2.192.73"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.672"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.473"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1322"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.223"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.365"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.360"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.098"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05153"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.270"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7445"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.784"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01941"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.801"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "77.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.407"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4455"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.958"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8316"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.723"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.445"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "983.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06022"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.57%  "
